# Peru Liga 1 - atualização de bases (29-02-2024)
#
# 1) A handful of existing data rows had their B:AC content (everything
#    except the running index in column A) swapped/rotated with a sibling
#    row - this mirrors upstream re-sorting the raw odds feed by match id.
# 2) 9 brand-new fixtures (ids 234-242, rows 236-244) are appended at the
#    bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($sheet, [int]$row1, [int]$row2)
    $r1 = $sheet.Range("B$row1`:AC$row1")
    $r2 = $sheet.Range("B$row2`:AC$row2")
    $v1 = $r1.Value()
    $v2 = $r2.Value()
    $r1.Value = $v2
    $r2.Value = $v1
}

# --- simple pairwise swaps -------------------------------------------------
Swap-RowData $ws 61 62
Swap-RowData $ws 156 157
Swap-RowData $ws 175 177
Swap-RowData $ws 181 182
Swap-RowData $ws 228 229

# --- 4-row rotation among 184,185,186,187 ----------------------------------
# new(184) = old(185); new(185) = old(187); new(187) = old(186); new(186) = old(184)
$v184 = $ws.Range("B184:AC184").Value()
$v185 = $ws.Range("B185:AC185").Value()
$v186 = $ws.Range("B186:AC186").Value()
$v187 = $ws.Range("B187:AC187").Value()

$ws.Range("B184:AC184").Value = $v185
$ws.Range("B185:AC185").Value = $v187
$ws.Range("B186:AC186").Value = $v184
$ws.Range("B187:AC187").Value = $v186

# --- append new fixtures (rows 236-244) ------------------------------------
function Add-Fixture {
    param($sheet, [int]$row, [int]$idx, [int]$matchId, [string]$home, [string]$away,
          [double]$dateSerial,
          [double]$K, [double]$L, [double]$M, [double]$N, [double]$O, [double]$P,
          [double]$Q, [double]$R, [double]$S, [double]$T, [double]$U, [double]$V)

    $sheet.Cells.Item($row, 1).Value = $idx
    $sheet.Cells.Item($row, 2).Value = $matchId
    $sheet.Cells.Item($row, 3).Value = "Peru Liga 1"
    $sheet.Cells.Item($row, 4).Value = "Peru Liga 1"
    $sheet.Cells.Item($row, 5).Value = $dateSerial
    $sheet.Cells.Item($row, 6).Value = $home
    $sheet.Cells.Item($row, 7).Value = $away

    # carry over the same cell formatting used by the existing data rows
    # (bold/centered index in column A, date number-format in column E)
    $sheet.Range("A235").Copy()
    $sheet.Cells.Item($row, 1).PasteSpecial(-4122)
    $sheet.Range("E235").Copy()
    $sheet.Cells.Item($row, 5).PasteSpecial(-4122)
    $sheet.Application.CutCopyMode = $false

    $sheet.Cells.Item($row, 11).Value = $K
    $sheet.Cells.Item($row, 12).Value = $L
    $sheet.Cells.Item($row, 13).Value = $M
    $sheet.Cells.Item($row, 14).Value = $N
    $sheet.Cells.Item($row, 15).Value = $O
    $sheet.Cells.Item($row, 16).Value = $P
    $sheet.Cells.Item($row, 17).Value = $Q
    $sheet.Cells.Item($row, 18).Value = $R
    $sheet.Cells.Item($row, 19).Value = $S
    $sheet.Cells.Item($row, 20).Value = $T
    $sheet.Cells.Item($row, 21).Value = $U
    $sheet.Cells.Item($row, 22).Value = $V
    $sheet.Cells.Item($row, 23).Value = 0
    $sheet.Cells.Item($row, 24).Value = 0
    $sheet.Cells.Item($row, 25).Value = 0
    $sheet.Cells.Item($row, 26).Value = 0
    $sheet.Cells.Item($row, 27).Value = 0
}

Add-Fixture $ws 236 234 7850288 "Alianza Atletico"           "Deportivo Garcilaso" 45351.67708333334 1.8   3.5 4.333 1.95 3.25 4    -0.5  2     1.85  2.5  2.05  1.8
Add-Fixture $ws 237 235 7850289 "AD Tarma"                   "Alianza Lima"        45352.70833333334 3.25  3.6 2     3.1  3.4  2.1  0.25  1.95  1.9   2.25 1.975 1.875
Add-Fixture $ws 238 236 7850290 "Universitario de Deportes"  "Sport Huancayo"      45352.9375         1.4   4   9     1.5  3.75 7    -1    1.85  2     2.5  1.95  1.9
Add-Fixture $ws 239 237 7850291 "Sporting Cristal"           "Atletico Grau"       45353.63541666666 1.25  5.5 11    1.25 6    10   -1.75 1.975 1.875 3    2     1.85
Add-Fixture $ws 240 238 7850292 "Comerciantes Unidos"        "FBC Melgar"          45353.73958333334 2.9   3.2 2.375 3.4  3.2  2.15 0.25  1.975 1.875 2.5  2     1.85
Add-Fixture $ws 241 239 7850293 "Cesar Vallejo"              "Cusco FC"            45353.875          1.8   3.75 4    1.615 4    4.75 -0.75 1.8   2.05  2.75 2.025 1.825
Add-Fixture $ws 242 240 7850294 "CD Los Chankas"             "Sport Boys"          45354.70833333334 1.909 3.5 3.75  1.571 3.8  5.5  -1    2.025 1.825 2.75 1.975 1.875
Add-Fixture $ws 243 241 7850295 "Cienciano"                  "Carlos Manucci"      45354.875          1.3   5   10    1.25 5.25 12   -1.75 2.025 1.825 3    2     1.85
Add-Fixture $ws 244 242 7850296 "Union Comercio"             "UTC Cajamarca"       45355.66666666666 2.5   3.2 2.75  2.15 3.3  3.2  -0.25 1.9   1.95  2.5  1.975 1.875
